# Add profile navigation label to English and Vietnamese translations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

# Insert a new row for the "nav.profile" translation key right before the
# "hero.title" row (currently row 13), shifting subsequent rows down by one.
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "nav.profile"
$ws.Range("B13").Value = "Profile"
$ws.Range("C13").Value = "Tài khoản của bạn"

# Fix the Vietnamese translation for "nav.library" (row 5, column C) which was
# previously left as the English value "Library".
$ws.Range("C5").Value = "Thư viện"

# Update the active selection as recorded by the author after the edit.
$ws.Range("C6").Select()
